# Excel to enum + Making sure the main scenes are selected when build
#
# 1) Both "Entities" and "Jp" sheets get a new row 2 inserted: id=100000,
#    name="None" (a new enum sentinel value). The former row 2 (the real
#    "Public_Phone_Maintenance_Book" entity) shifts down to row 3 and its
#    id becomes 100001.
# 2) The active sheet/tab moves from "Entities" to "Jp" (so the Japanese
#    scene is the one selected when the project builds), and the stale
#    selection on "Entities" is moved onto the now-shifted data row.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "Entities"
$ws2 = $wb.Worksheets.Item(2)   # "Jp"

foreach ($ws in @($ws1, $ws2)) {
    # Push the existing row 2 down to row 3 (carries its styles/row height
    # along with it), leaving a blank row 2 behind.
    $ws.Rows.Item(2).Insert() | Out-Null

    # New row 2: the "None" enum entry.
    $ws.Range("A2").Style = "Good"
    $ws.Range("B2").Style = "Good"
    $ws.Range("A2").Value = 100000
    $ws.Range("B2").Value = "None"

    # The entity that used to be row 2 is now row 3; bump its id.
    $ws.Range("A3").Value = 100001
}

# sheet1 ("Entities") originally spanned columns A:L on its data row, so
# the shifted row 3 needs the (empty, but styled) trailing cells restored.
$ws1.Range("G3:L3").Style = "Good"

# Fix up the selections left stale on sheet1 after the insert, then make
# "Jp" the active tab (matches the committed workbookView/sheetView state).
$ws1.Range("D3").Select() | Out-Null
$ws2.Activate() | Out-Null
